$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.085.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.516.10'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.29%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '489.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +12.75%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +6.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.535.88'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0986'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.64'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.335'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.92%  '
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.952.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.143.79'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +8.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000137'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +8.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.527.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.40%  '
$ws.Range("E19").Value = '  +6.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.99%  '
$ws.Range("E22").Value = '  -1.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.68'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.413'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.166'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +10.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.619.76'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.58%  '
$ws.Range("E29").Value = '  +6.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0795'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +12.82%  '
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.55'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.84%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.50'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.33%  '
$ws.Range("E34").Value = '  +8.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.71%  '
$ws.Range("E36").Value = '  +12.54%  '
$ws.Range("E37").Value = '  +12.23%  '
$ws.Range("E38").Value = '  +7.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.39'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.55'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.621'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0558'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.995'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("E44").Value = '  +10.14%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +14.35%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '266.70'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +29.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0911'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.68%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0227'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.957.25'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.18%  '
